# Aula 36 - Inserindo departamentos
#
# Appends 5 new rows (26-30) to the "Planilha1" log sheet, all belonging to
# the new session "8. Departamento: Controller & View" / lesson
# "36. Inserindo departamentos", each holding one timestamped note about
# th:object / th:field usage with Thymeleaf.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Seed formatting for the new rows by copying it from row 14, which
#    already uses the "section header" styling (bold red font) that the
#    new rows need: columns B/C without wrap, columns D/E with wrap.
#    Using Copy + PasteSpecial(formats) re-uses the existing style/font
#    entries instead of creating new ones.
# ---------------------------------------------------------------------
$ws.Range("B14:C14").Copy()
$ws.Range("B26:C26").PasteSpecial(-4122)
$ws.Range("B27:C27").PasteSpecial(-4122)
$ws.Range("B28:C28").PasteSpecial(-4122)
$ws.Range("B29:C29").PasteSpecial(-4122)
$ws.Range("B30:C30").PasteSpecial(-4122)

$ws.Range("E14").Copy()
$ws.Range("D26:E26").PasteSpecial(-4122)
$ws.Range("D27:E27").PasteSpecial(-4122)
$ws.Range("D28:E28").PasteSpecial(-4122)
$ws.Range("D29:E29").PasteSpecial(-4122)
$ws.Range("D30:E30").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) Fill in the data for the 5 new rows.
# ---------------------------------------------------------------------
$sessao = "8. Departamento: Controller & View"
$aula = "36. Inserindo departamentos"

$obs1 = "
2:12 - aplicação de chamada de endpoint atraves de tags e expressões do Thymeleaf. a expressão utilizada para o thymeleaf é ""@{/seu/endpoit/aqui}"" ... o parêntese permanece no código."

$obs2 = "3:32
declaração de váriavel do thymeleaf no documento HTML com a expressão: th:object=""`${departamento}"". Essa variável podera ser acessada pelo backend, como por exemplo acessar objetos da entidade. ela tem ligação direta com outra expressão do thymeleaf chamada th:field=""*{nomeDoAtributoDoObjetoAqui}"""

$obs3 = "4:30
declaração de atributo do thymeleaf no documento HTML, diferente de variavel, o th:field=""*{nomeDoAtributoDoObjetoAqui}""tem ligação direta com o th:object. Com o th:field é possível acessar os atributos de uma classe."

$obs4 = "4:37
em resumo: ""th:object"" e o ""th:field"" são capazes de instanciar um objeto de entidade/banco de dados e setar os seus atributos de acordo com os valores dos elementos HTML e entregar os dados para o Controller da aplicação."

$obs5 = "6:58
uma excessão que vale a pena comentar, pois é sobre o fato de ter declarado um th:object na pagina HTML de cadastro de Departamentos porém não foi declarado no endpoint, ou seja, ao chamar o endpoint, gera este erro. A solução é declarar tbm em forma de parametro no metodo do endpoint."

foreach ($r in 26..30) {
    $ws.Range("B$r").Value = 36
    $ws.Range("C$r").Value = $sessao
    $ws.Range("D$r").Value = $aula
}

$ws.Range("E26").Value = $obs1
$ws.Range("E27").Value = $obs2
$ws.Range("E28").Value = $obs3
$ws.Range("E29").Value = $obs4
$ws.Range("E30").Value = $obs5

# ---------------------------------------------------------------------
# 3) Row heights (auto-calculated by Excel for the wrapped text in the
#    observation column, based on how many lines each note spans).
# ---------------------------------------------------------------------
$ws.Rows.Item(26).RowHeight = 60
$ws.Rows.Item(27).RowHeight = 75
$ws.Rows.Item(28).RowHeight = 60
$ws.Rows.Item(29).RowHeight = 60
$ws.Rows.Item(30).RowHeight = 75

# ---------------------------------------------------------------------
# 4) Leave the view scrolled down to the newly added rows, with D27 as
#    the active cell (matching where the author ended up editing).
# ---------------------------------------------------------------------
$null = $ws.Range("D27").Select()
$excel.ActiveWindow.ScrollRow = 24
$excel.ActiveWindow.ScrollColumn = 1
